$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the existing
# header style (bold/bordered/centered) used by A1:H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the I0 and IF data columns for rows 2-69.
$iVals = @(8,9,8,7,7,5,8,9,7,4,9,8,9,8,7,6,7,8,9,4,8,7,7,8,9,7,9,7,9,10,8,9,7,9,7,8,8,9,11,10,9,9,7,9,5,4,6,8,9,5,8,9,7,6,7,8,7,6,8,7,8,7,8,5,8,8,7,4)
$jVals = @(8,9,8,7,7,6,8,9,7,5,9,8,9,8,8,6,7,8,9,5,8,7,7,8,9,7,9,7,9,10,8,9,7,9,8,8,9,9,11,10,9,9,7,9,5,4,6,8,9,5,8,9,7,6,7,8,7,6,8,7,8,7,8,5,8,8,7,4)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
